# Deploy the implementation guide.
# Updates the "Metadata" sheet of the CodeSystem-disease-status workbook:
#   - Date value regenerated
#   - Contact value regenerated (now resolves a display name + URL)
#   - a new "Jurisdiction" property row is inserted right after "Contact"
#     (shifting Description..Count down by one row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Date (row 8, column B) gets a new timestamp.
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# 2) Contact (row 10, column B) now has a resolved display value.
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 3) Insert a new row above the old row 11 ("Description") for "Jurisdiction",
#    pushing Description..Count down to rows 12..22.
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Match the surrounding "body row" look (wrap text, top-aligned, bordered)
# used by every other data row on this sheet.
$bodyRow = $ws.Range("A10:B10")
$newRow = $ws.Range("A11:B11")
$newRow.WrapText = $true
$newRow.VerticalAlignment = -4160
$newRow.Borders.LineStyle = 1
